$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "81-45="
$tbl.Cell(1,2).Range.Text = "77+8="
$tbl.Cell(1,3).Range.Text = "42-7="
$tbl.Cell(1,4).Range.Text = "46-29="
$tbl.Cell(1,5).Range.Text = "73+19="
$tbl.Cell(2,1).Range.Text = "56-18="
$tbl.Cell(2,2).Range.Text = "83-75="
$tbl.Cell(2,3).Range.Text = "88-69="
$tbl.Cell(2,4).Range.Text = "86-9="
$tbl.Cell(2,5).Range.Text = "40-3="
$tbl.Cell(3,1).Range.Text = "8+26="
$tbl.Cell(3,2).Range.Text = "71-69="
$tbl.Cell(3,3).Range.Text = "29+67="
$tbl.Cell(3,4).Range.Text = "19+14="
$tbl.Cell(3,5).Range.Text = "72-69="
$tbl.Cell(4,1).Range.Text = "36+6="
$tbl.Cell(4,2).Range.Text = "80-29="
$tbl.Cell(4,3).Range.Text = "70-41="
$tbl.Cell(4,4).Range.Text = "61-53="
$tbl.Cell(4,5).Range.Text = "9+5="
$tbl.Cell(5,1).Range.Text = "16+7="
$tbl.Cell(5,2).Range.Text = "50-11="
$tbl.Cell(5,3).Range.Text = "56-8="
$tbl.Cell(5,4).Range.Text = "79+3="
$tbl.Cell(5,5).Range.Text = "90-43="
$tbl.Cell(6,1).Range.Text = "70-37="
$tbl.Cell(6,2).Range.Text = "84-48="
$tbl.Cell(6,3).Range.Text = "80-5="
$tbl.Cell(6,4).Range.Text = "53-4="
$tbl.Cell(6,5).Range.Text = "91-89="
$tbl.Cell(7,1).Range.Text = "90-14="
$tbl.Cell(7,2).Range.Text = "82-39="
$tbl.Cell(7,3).Range.Text = "31-14="
$tbl.Cell(7,4).Range.Text = "65-28="
$tbl.Cell(7,5).Range.Text = "37-28="
$tbl.Cell(8,1).Range.Text = "81-57="
$tbl.Cell(8,2).Range.Text = "93-57="
$tbl.Cell(8,3).Range.Text = "65+27="
$tbl.Cell(8,4).Range.Text = "92-83="
$tbl.Cell(8,5).Range.Text = "98-89="
$tbl.Cell(9,1).Range.Text = "30-17="
$tbl.Cell(9,2).Range.Text = "38+24="
$tbl.Cell(9,3).Range.Text = "91-62="
$tbl.Cell(9,4).Range.Text = "77+14="
$tbl.Cell(9,5).Range.Text = "46+39="
$tbl.Cell(10,1).Range.Text = "63-56="
$tbl.Cell(10,2).Range.Text = "15+69="
$tbl.Cell(10,3).Range.Text = "5+37="
$tbl.Cell(10,4).Range.Text = "26+26="
$tbl.Cell(10,5).Range.Text = "51-3="
$tbl.Cell(11,1).Range.Text = "75-6="
$tbl.Cell(11,2).Range.Text = "9+68="
$tbl.Cell(11,3).Range.Text = "47+37="
$tbl.Cell(11,4).Range.Text = "91-82="
$tbl.Cell(11,5).Range.Text = "7+18="
$tbl.Cell(12,1).Range.Text = "76+9="
$tbl.Cell(12,2).Range.Text = "54-27="
$tbl.Cell(12,3).Range.Text = "77+16="
$tbl.Cell(12,4).Range.Text = "93-88="
$tbl.Cell(12,5).Range.Text = "80-49="
$tbl.Cell(13,1).Range.Text = "28+8="
$tbl.Cell(13,2).Range.Text = "47+49="
$tbl.Cell(13,3).Range.Text = "14-6="
$tbl.Cell(13,4).Range.Text = "37+46="
$tbl.Cell(13,5).Range.Text = "49+34="
$tbl.Cell(14,1).Range.Text = "90-51="
$tbl.Cell(14,2).Range.Text = "65-58="
$tbl.Cell(14,3).Range.Text = "16+48="
$tbl.Cell(14,4).Range.Text = "86+6="
$tbl.Cell(14,5).Range.Text = "51-5="
$tbl.Cell(15,1).Range.Text = "19+77="
$tbl.Cell(15,2).Range.Text = "34-19="
$tbl.Cell(15,3).Range.Text = "29+67="
$tbl.Cell(15,4).Range.Text = "14+19="
$tbl.Cell(15,5).Range.Text = "78+6="
$tbl.Cell(16,1).Range.Text = "39+59="
$tbl.Cell(16,2).Range.Text = "46+8="
$tbl.Cell(16,3).Range.Text = "80-11="
$tbl.Cell(16,4).Range.Text = "48+8="
$tbl.Cell(16,5).Range.Text = "76+6="
$tbl.Cell(17,1).Range.Text = "79+8="
$tbl.Cell(17,2).Range.Text = "88+7="
$tbl.Cell(17,3).Range.Text = "80-35="
$tbl.Cell(17,4).Range.Text = "42-33="
$tbl.Cell(17,5).Range.Text = "6+77="
$tbl.Cell(18,1).Range.Text = "49+27="
$tbl.Cell(18,2).Range.Text = "18+17="
$tbl.Cell(18,3).Range.Text = "73-36="
$tbl.Cell(18,4).Range.Text = "62-33="
$tbl.Cell(18,5).Range.Text = "89+8="
$tbl.Cell(19,1).Range.Text = "83-37="
$tbl.Cell(19,2).Range.Text = "68+19="
$tbl.Cell(19,3).Range.Text = "68-49="
$tbl.Cell(19,4).Range.Text = "38-19="
$tbl.Cell(19,5).Range.Text = "66-17="
$tbl.Cell(20,1).Range.Text = "25+69="
$tbl.Cell(20,2).Range.Text = "15+17="
$tbl.Cell(20,3).Range.Text = "71-69="
$tbl.Cell(20,4).Range.Text = "80-47="
$tbl.Cell(20,5).Range.Text = "40-29="